$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.6618222594261169
$ws.Range("B1").Value = 1.32581102848053
$ws.Range("C1").Value = 4.058460235595703
$ws.Range("D1").Value = 1.761385798454285
$ws.Range("E1").Value = 0.485239565372467
